$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "43.195.53"
$ws.Range("E2").Value = "  -1.20%  "
$ws.Range("D3").Value = "2.279.59"
$ws.Range("E3").Value = "  -0.59%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "112.84"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.25%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "265.74"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.21%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.616"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.606"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -1.26%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "47.48"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.01%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0934"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.37%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "9.30"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +8.61%  "
$ws.Range("E13").Value = "  +1.44%  "
$ws.Range("E14").Value = "  -0.68%  "
$ws.Range("D15").Value = "2.610.18"
$ws.Range("E15").Value = "  -0.97%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.867"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +2.17%  "
$ws.Range("D17").Value = "2.274.92"
$ws.Range("E17").Value = "  -0.54%  "
$ws.Range("D18").Value = "43.145.21"
$ws.Range("E18").Value = "  -1.06%  "
$ws.Range("E19").Value = "  -1.21%  "
$ws.Range("E20").Value = "  +3.81%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.85"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.77%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.49"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "234.27"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.38%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.60"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.23%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.88"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.29%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.02"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +1.58%  "
$ws.Range("E27").Value = "  -0.79%  "
$ws.Range("E28").Value = "  +0.64%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "41.03"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -3.01%  "
$ws.Range("E30").Value = "  -2.30%  "
$ws.Range("E31").Value = "  -0.83%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "174.15"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -1.26%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "21.47"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.81%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0905"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -2.06%  "
$ws.Range("E35").Value = "  +3.49%  "
$ws.Range("E36").Value = "  +0.28%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.66"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -1.15%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0365"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +2.94%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.90"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +2.52%  "
$ws.Range("E40").Value = "  -5.19%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.63"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +8.49%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "76.76"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +4.21%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "14.16"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +3.62%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.238"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -1.99%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.21"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +4.47%  "
$ws.Range("E46").Value = "  -0.12%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.39"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -2.99%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.67"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -1.09%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "103.62"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.35%  "
$ws.Range("E50").Value = "  +1.23%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0997"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.67%  "
